$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsARM = $wb.Worksheets.Item("ARM")
$wsBSM = $wb.Worksheets.Item("BSM")
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCUL = $wb.Worksheets.Item("CUL")
$wsGSM = $wb.Worksheets.Item("GSM")
$wsLTW = $wb.Worksheets.Item("LTW")
$wsWVR = $wb.Worksheets.Item("WVR")

# ALC row 9
$wsALC.Range("H9").Value = 130.81818
$wsALC.Range("I9").Value = 123.9
$wsALC.Range("J9").Value = 200
$wsALC.Range("K9").Value = 123.9
$wsALC.Range("L9").Value = 200
$wsALC.Range("M9").Value = 45.09999999999999
$wsALC.Range("N9").Value = -538

# ALC row 100
$wsALC.Range("H100").Value = 3133.75
$wsALC.Range("I100").Value = 2846.111
$wsALC.Range("J100").Value = 3996.6667
$wsALC.Range("K100").Value = 2846.111
$wsALC.Range("L100").Value = 3996.6667
$wsALC.Range("M100").Value = -2305.111
$wsALC.Range("N100").Value = -5078.6667

# ALC row 111
$wsALC.Range("H111").Value = 13278
$wsALC.Range("I111").Value = 3571.6
$wsALC.Range("J111").Value = 20211.143
$wsALC.Range("K111").Value = 10714.8
$wsALC.Range("L111").Value = 60633.429
$wsALC.Range("M111").Value = -7647.799999999999
$wsALC.Range("N111").Value = -66767.429

# ALC row 137
$wsALC.Range("H137").Value = 2774.1562
$wsALC.Range("I137").Value = 1693.4615
$wsALC.Range("J137").Value = 3513.5789
$wsALC.Range("K137").Value = 5080.3845
$wsALC.Range("L137").Value = 10540.7367
$wsALC.Range("M137").Value = -2530.3845
$wsALC.Range("N137").Value = -15640.7367

# ALC row 138
$wsALC.Range("H138").Value = 3031.63
$wsALC.Range("I138").Value = 1941.1154
$wsALC.Range("J138").Value = 3414.7837
$wsALC.Range("K138").Value = 5823.3462
$wsALC.Range("L138").Value = 10244.3511
$wsALC.Range("M138").Value = -683.3462
$wsALC.Range("N138").Value = -20524.3511

# ARM row 61
$wsARM.Range("H61").Value = 3812.4443
$wsARM.Range("I61").Value = 4052
$wsARM.Range("J61").Value = 3333.3333
$wsARM.Range("K61").Value = 4052
$wsARM.Range("L61").Value = 3333.3333
$wsARM.Range("M61").Value = -3840
$wsARM.Range("N61").Value = -3757.3333

# ARM row 74
$wsARM.Range("H74").Value = 1428.68
$wsARM.Range("I74").Value = 1311.3334
$wsARM.Range("J74").Value = 1844.7273
$wsARM.Range("K74").Value = 1311.3334
$wsARM.Range("L74").Value = 1844.7273
$wsARM.Range("M74").Value = -437.3334
$wsARM.Range("N74").Value = -3592.7273

# ARM row 77
$wsARM.Range("H77").Value = 1428.68
$wsARM.Range("I77").Value = 1311.3334
$wsARM.Range("J77").Value = 1844.7273
$wsARM.Range("K77").Value = 6556.666999999999
$wsARM.Range("L77").Value = 9223.636500000001
$wsARM.Range("M77").Value = -2188.666999999999
$wsARM.Range("N77").Value = -17959.6365

# ARM row 132
$wsARM.Range("H132").Value = 5123.447
$wsARM.Range("I132").Value = 5896.8
$wsARM.Range("J132").Value = 3758.7058
$wsARM.Range("K132").Value = 17690.4
$wsARM.Range("L132").Value = 11276.1174
$wsARM.Range("M132").Value = -15160.4
$wsARM.Range("N132").Value = -16336.1174

# ARM row 136
$wsARM.Range("H136").Value = 3812.4443
$wsARM.Range("I136").Value = 4052
$wsARM.Range("J136").Value = 3333.3333
$wsARM.Range("K136").Value = 12156
$wsARM.Range("L136").Value = 9999.999899999999
$wsARM.Range("M136").Value = -9606
$wsARM.Range("N136").Value = -15099.9999

# BSM row 20
$wsBSM.Range("H20").Value = 84316.664
$wsBSM.Range("I20").Value = 200799.4
$wsBSM.Range("J20").Value = 1114.7142
$wsBSM.Range("K20").Value = 200799.4
$wsBSM.Range("L20").Value = 1114.7142
$wsBSM.Range("M20").Value = -200552.4
$wsBSM.Range("N20").Value = -1608.7142

# BSM row 105
$wsBSM.Range("H105").Value = 5104923
$wsBSM.Range("I105").Value = 6496183
$wsBSM.Range("K105").Value = 6496183
$wsBSM.Range("M105").Value = -6494436

# BSM row 132
$wsBSM.Range("H132").Value = 0
$wsBSM.Range("J132").Value = 0
$wsBSM.Range("L132").Value = 0
$wsBSM.Range("N132").ClearContents()

# BSM row 135
$wsBSM.Range("H135").Value = 73044.664
$wsBSM.Range("J135").Value = 73044.664
$wsBSM.Range("L135").Value = 73044.664
$wsBSM.Range("N135").Value = -83184.664

# BSM row 137
$wsBSM.Range("H137").Value = 74780
$wsBSM.Range("J137").Value = 74780
$wsBSM.Range("L137").Value = 74780
$wsBSM.Range("N137").Value = -84980

# BSM row 140
$wsBSM.Range("H140").Value = 121260
$wsBSM.Range("J140").Value = 121260
$wsBSM.Range("L140").Value = 121260
$wsBSM.Range("N140").Value = -131620

# CRP row 58
$wsCRP.Range("H58").Value = 1765782.6
$wsCRP.Range("I58").Value = 2316889.8
$wsCRP.Range("J58").Value = 2240
$wsCRP.Range("K58").Value = 2316889.8
$wsCRP.Range("L58").Value = 2240
$wsCRP.Range("M58").Value = -2316686.8
$wsCRP.Range("N58").Value = -2646

# CRP row 99
$wsCRP.Range("H99").Value = 2022.5
$wsCRP.Range("I99").Value = 2114.2856
$wsCRP.Range("J99").Value = 1380
$wsCRP.Range("K99").Value = 2114.2856
$wsCRP.Range("L99").Value = 1380
$wsCRP.Range("M99").Value = -616.2856000000002
$wsCRP.Range("N99").Value = -4376

# CRP row 126
$wsCRP.Range("H126").Value = 2022.5
$wsCRP.Range("I126").Value = 2114.2856
$wsCRP.Range("J126").Value = 1380
$wsCRP.Range("K126").Value = 6342.8568
$wsCRP.Range("L126").Value = 4140
$wsCRP.Range("M126").Value = -3872.8568
$wsCRP.Range("N126").Value = -9080

# CRP row 134
$wsCRP.Range("H134").Value = 2189.1333
$wsCRP.Range("I134").Value = 1840.2307
$wsCRP.Range("J134").Value = 4457
$wsCRP.Range("K134").Value = 5520.6921
$wsCRP.Range("L134").Value = 13371
$wsCRP.Range("M134").Value = -2985.6921
$wsCRP.Range("N134").Value = -18441

# CRP row 136
$wsCRP.Range("H136").Value = 1765782.6
$wsCRP.Range("I136").Value = 2316889.8
$wsCRP.Range("J136").Value = 2240
$wsCRP.Range("K136").Value = 6950669.399999999
$wsCRP.Range("L136").Value = 6720
$wsCRP.Range("M136").Value = -6948119.399999999
$wsCRP.Range("N136").Value = -11820

# CUL row 5
$wsCUL.Range("H5").Value = 1090.6111
$wsCUL.Range("I5").Value = 1493.45
$wsCUL.Range("J5").Value = 587.0625
$wsCUL.Range("K5").Value = 4480.35
$wsCUL.Range("L5").Value = 1761.1875
$wsCUL.Range("M5").Value = -4368.35
$wsCUL.Range("N5").Value = -1985.1875

# CUL row 80
$wsCUL.Range("H80").Value = 13960
$wsCUL.Range("I80").Value = 12100
$wsCUL.Range("J80").Value = 16750
$wsCUL.Range("K80").Value = 36300
$wsCUL.Range("L80").Value = 50250
$wsCUL.Range("M80").Value = -35364
$wsCUL.Range("N80").Value = -52122

# CUL row 83
$wsCUL.Range("H83").Value = 13960
$wsCUL.Range("I83").Value = 12100
$wsCUL.Range("J83").Value = 16750
$wsCUL.Range("K83").Value = 108900
$wsCUL.Range("L83").Value = 150750
$wsCUL.Range("M83").Value = -104220
$wsCUL.Range("N83").Value = -160110

# CUL row 86
$wsCUL.Range("H86").Value = 35866.668
$wsCUL.Range("I86").Value = 3799.5
$wsCUL.Range("J86").Value = 100001
$wsCUL.Range("K86").Value = 11398.5
$wsCUL.Range("L86").Value = 300003
$wsCUL.Range("M86").Value = -10212.5
$wsCUL.Range("N86").Value = -302375

# CUL row 89
$wsCUL.Range("H89").Value = 35866.668
$wsCUL.Range("I89").Value = 3799.5
$wsCUL.Range("J89").Value = 100001
$wsCUL.Range("K89").Value = 34195.5
$wsCUL.Range("L89").Value = 900009
$wsCUL.Range("M89").Value = -28267.5
$wsCUL.Range("N89").Value = -911865

# CUL row 107
$wsCUL.Range("H107").Value = 1330
$wsCUL.Range("I107").Value = 1091.4468
$wsCUL.Range("J107").Value = 1730.4286
$wsCUL.Range("K107").Value = 3274.3404
$wsCUL.Range("L107").Value = 5191.2858
$wsCUL.Range("M107").Value = -1354.3404
$wsCUL.Range("N107").Value = -9031.2858

# CUL row 112
$wsCUL.Range("H112").Value = 4985
$wsCUL.Range("I112").Value = 2980
$wsCUL.Range("J112").Value = 8326.666999999999
$wsCUL.Range("K112").Value = 8940
$wsCUL.Range("L112").Value = 24980.001
$wsCUL.Range("M112").Value = -7832
$wsCUL.Range("N112").Value = -27196.001

# CUL row 132
$wsCUL.Range("H132").Value = 1284.2433
$wsCUL.Range("I132").Value = 1133.8334
$wsCUL.Range("J132").Value = 1313.3549
$wsCUL.Range("K132").Value = 10204.5006
$wsCUL.Range("L132").Value = 11820.1941
$wsCUL.Range("M132").Value = -7674.500599999999
$wsCUL.Range("N132").Value = -16880.1941

# CUL row 135
$wsCUL.Range("H135").Value = 1090.6111
$wsCUL.Range("I135").Value = 1493.45
$wsCUL.Range("J135").Value = 587.0625
$wsCUL.Range("K135").Value = 13441.05
$wsCUL.Range("L135").Value = 5283.5625
$wsCUL.Range("M135").Value = -10906.05
$wsCUL.Range("N135").Value = -10353.5625

# GSM row 123
$wsGSM.Range("H123").Value = 8598.162
$wsGSM.Range("J123").Value = 8598.162
$wsGSM.Range("L123").Value = 8598.162
$wsGSM.Range("N123").Value = -13498.162

# GSM row 132
$wsGSM.Range("H132").Value = 3368.2104
$wsGSM.Range("I132").Value = 2502.6667
$wsGSM.Range("J132").Value = 4147.2
$wsGSM.Range("K132").Value = 7508.000100000001
$wsGSM.Range("L132").Value = 12441.6
$wsGSM.Range("M132").Value = -4978.000100000001
$wsGSM.Range("N132").Value = -17501.6

# LTW row 7
$wsLTW.Range("H7").Value = 3483.3333
$wsLTW.Range("J7").Value = 3483.3333
$wsLTW.Range("L7").Value = 3483.3333
$wsLTW.Range("N7").Value = -3707.3333

# LTW row 33
$wsLTW.Range("H33").Value = 50015
$wsLTW.Range("I33").Value = 50015
$wsLTW.Range("K33").Value = 50015
$wsLTW.Range("M33").Value = -49725

# LTW row 40
$wsLTW.Range("H40").Value = 3066.5557
$wsLTW.Range("I40").Value = 2720
$wsLTW.Range("J40").Value = 3499.75
$wsLTW.Range("K40").Value = 2720
$wsLTW.Range("L40").Value = 3499.75
$wsLTW.Range("M40").Value = -2584
$wsLTW.Range("N40").Value = -3771.75

# LTW row 126
$wsLTW.Range("H126").Value = 3483.3333
$wsLTW.Range("J126").Value = 3483.3333
$wsLTW.Range("L126").Value = 10449.9999
$wsLTW.Range("N126").Value = -15389.9999

# LTW row 136
$wsLTW.Range("H136").Value = 4571.4287
$wsLTW.Range("I136").Value = 0
$wsLTW.Range("K136").Value = 0
$wsLTW.Range("M136").ClearContents()

# WVR row 64
$wsWVR.Range("H64").Value = 30000
$wsWVR.Range("J64").Value = 30000
$wsWVR.Range("L64").Value = 30000
$wsWVR.Range("N64").Value = -30496

# WVR row 67
$wsWVR.Range("H67").Value = 30000
$wsWVR.Range("J67").Value = 30000
$wsWVR.Range("L67").Value = 30000
$wsWVR.Range("N67").Value = -31716

# WVR row 126
$wsWVR.Range("H126").Value = 8316
$wsWVR.Range("I126").Value = 11475
$wsWVR.Range("J126").Value = 1998
$wsWVR.Range("K126").Value = 34425
$wsWVR.Range("L126").Value = 5994
$wsWVR.Range("M126").Value = -31955
$wsWVR.Range("N126").Value = -10934

# WVR row 132
$wsWVR.Range("H132").Value = 3084.7896
$wsWVR.Range("I132").Value = 2841.5186
$wsWVR.Range("J132").Value = 3681.9092
$wsWVR.Range("K132").Value = 8524.5558
$wsWVR.Range("L132").Value = 11045.7276
$wsWVR.Range("M132").Value = -5994.5558
$wsWVR.Range("N132").Value = -16105.7276

# WVR row 136
$wsWVR.Range("H136").Value = 4299.2354
$wsWVR.Range("I136").Value = 5512.4287
$wsWVR.Range("K136").Value = 16537.2861
$wsWVR.Range("M136").Value = -13987.2861

